$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 95
$ws.Range("I2").Value = 95
$ws.Range("K2").Value = 95
$ws.Range("M2").Value = 18

$ws.Range("H38").Value = 2541.2856
$ws.Range("I38").Value = 175.33333
$ws.Range("J38").Value = 6800
$ws.Range("K38").Value = 525.99999
$ws.Range("L38").Value = 20400
$ws.Range("M38").Value = -153.99999
$ws.Range("N38").Value = -21144

$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("N43").Value = 0
$ws.Range("L43").Value = $null

$ws.Range("H58").Value = 3230
$ws.Range("I58").Value = 1476.6666
$ws.Range("J58").Value = 4983.3335
$ws.Range("K58").Value = 4429.9998
$ws.Range("L58").Value = 14950.0005
$ws.Range("M58").Value = -4279.9998
$ws.Range("N58").Value = -15250.0005

$ws.Range("H86").Value = 4000
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").Value = $null

$ws.Range("H89").Value = 4000
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").Value = $null

$ws.Range("H98").Value = 452.66666
$ws.Range("I98").Value = 452.66666
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 452.66666
$ws.Range("L98").Value = 0
$ws.Range("N98").Value = 1045.33334
$ws.Range("M98").Value = $null

$ws.Range("H106").Value = 2104
$ws.Range("I106").Value = 1855
$ws.Range("K106").Value = 1855
$ws.Range("M106").Value = -1224

$ws.Range("H115").Value = 1149
$ws.Range("I115").Value = 1149
$ws.Range("K115").Value = 3447
$ws.Range("M115").Value = -1880

$ws.Range("H122").Value = 452.66666
$ws.Range("I122").Value = 452.66666
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 1357.99998
$ws.Range("L122").Value = 0
$ws.Range("N122").Value = 1092.00002
$ws.Range("M122").Value = $null

$ws.Range("H125").Value = 1031.5
$ws.Range("I125").Value = 1031.5
$ws.Range("K125").Value = 9283.5
$ws.Range("M125").Value = -6823.5

$ws.Range("H137").Value = 2105.6667
$ws.Range("I137").Value = 629
$ws.Range("K137").Value = 1887
$ws.Range("M137").Value = 663

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 12722.8125
$ws.Range("J44").Value = 13171
$ws.Range("L44").Value = 13171
$ws.Range("N44").Value = -14147

$ws.Range("H74").Value = 7386.75
$ws.Range("I74").Value = 7541.5
$ws.Range("K74").Value = 7541.5
$ws.Range("M74").Value = -6667.5

$ws.Range("H77").Value = 7386.75
$ws.Range("I77").Value = 7541.5
$ws.Range("K77").Value = 37707.5
$ws.Range("M77").Value = -33339.5

$ws.Range("H97").Value = 771.5
$ws.Range("I97").Value = 771.5
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 771.5
$ws.Range("L97").Value = 0
$ws.Range("N97").Value = -275.5
$ws.Range("M97").Value = $null

$ws.Range("H132").Value = 1478.5714
$ws.Range("I132").Value = 1581.909
$ws.Range("J132").Value = 1099.6666
$ws.Range("K132").Value = 4745.727000000001
$ws.Range("L132").Value = 3298.9998
$ws.Range("M132").Value = -2215.727000000001
$ws.Range("N132").Value = -8358.9998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 20000
$ws.Range("I20").Value = 20000
$ws.Range("K20").Value = 20000
$ws.Range("M20").Value = -19753

$ws.Range("H86").Value = 3178.9285
$ws.Range("I86").Value = 2130.5
$ws.Range("K86").Value = 2130.5
$ws.Range("M86").Value = -1007.5

$ws.Range("H89").Value = 3178.9285
$ws.Range("I89").Value = 2130.5
$ws.Range("K89").Value = 10652.5
$ws.Range("M89").Value = -5036.5

$ws.Range("H102").Value = 35000
$ws.Range("J102").Value = 35000
$ws.Range("L102").Value = 35000
$ws.Range("N102").Value = -41490

$ws.Range("H105").Value = 6163873
$ws.Range("I105").Value = 10085065
$ws.Range("K105").Value = 10085065
$ws.Range("M105").Value = -10083318

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6006.9546
$ws.Range("I31").Value = 2230.8
$ws.Range("K31").Value = 2230.8
$ws.Range("M31").Value = -1935.8

$ws.Range("H34").Value = 6006.9546
$ws.Range("I34").Value = 2230.8
$ws.Range("K34").Value = 2230.8
$ws.Range("M34").Value = -2028.8

$ws.Range("H58").Value = 4763.1665
$ws.Range("I58").Value = 3296.2856
$ws.Range("K58").Value = 3296.2856
$ws.Range("M58").Value = -3093.2856

$ws.Range("H134").Value = 3572.75
$ws.Range("I134").Value = 3572.75
$ws.Range("K134").Value = 10718.25
$ws.Range("M134").Value = -8183.25

$ws.Range("H136").Value = 4763.1665
$ws.Range("I136").Value = 3296.2856
$ws.Range("K136").Value = 9888.856800000001
$ws.Range("M136").Value = -7338.856800000001

$ws.Range("H141").Value = 45895.4
$ws.Range("J141").Value = 46563
$ws.Range("L141").Value = 46563
$ws.Range("N141").Value = -56923

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 8014.4165
$ws.Range("J39").Value = 8615.727999999999
$ws.Range("L39").Value = 25847.184
$ws.Range("N39").Value = -26435.184

$ws.Range("H131").Value = 1848.826
$ws.Range("I131").Value = 1289.8
$ws.Range("J131").Value = 2004.1111
$ws.Range("K131").Value = 3869.4
$ws.Range("L131").Value = 6012.3333
$ws.Range("M131").Value = 1170.6
$ws.Range("N131").Value = -16092.3333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 153.97298
$ws.Range("I2").Value = 43.517242
$ws.Range("K2").Value = 43.517242
$ws.Range("M2").Value = 69.48275799999999

$ws.Range("H70").Value = 4835.6665
$ws.Range("I70").Value = 4835.6665
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 4835.6665
$ws.Range("L70").Value = 0
$ws.Range("N70").Value = -4565.6665
$ws.Range("M70").Value = $null

$ws.Range("H73").Value = 4835.6665
$ws.Range("I73").Value = 4835.6665
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 4835.6665
$ws.Range("L73").Value = 0
$ws.Range("N73").Value = -3899.6665
$ws.Range("M73").Value = $null

$ws.Range("H80").Value = 1615.4
$ws.Range("I80").Value = 1725.5
$ws.Range("K80").Value = 1725.5
$ws.Range("M80").Value = -727.5

$ws.Range("H83").Value = 1615.4
$ws.Range("I83").Value = 1725.5
$ws.Range("K83").Value = 8627.5
$ws.Range("M83").Value = -3635.5

$ws.Range("H99").Value = 2635.1428
$ws.Range("I99").Value = 2635.1428
$ws.Range("K99").Value = 2635.1428
$ws.Range("M99").Value = -389.1428000000001

$ws.Range("H122").Value = 2616
$ws.Range("I122").Value = 1947
$ws.Range("K122").Value = 5841
$ws.Range("M122").Value = -3391

$ws.Range("H132").Value = 2261.4546
$ws.Range("I132").Value = 2261.4546
$ws.Range("K132").Value = 6784.3638
$ws.Range("M132").Value = -4254.3638

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 933.3333
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 900
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 900
$ws.Range("M22").Value = -705
$ws.Range("N22").Value = -1490

$ws.Range("H27").Value = 933.3333
$ws.Range("I27").Value = 1000
$ws.Range("J27").Value = 900
$ws.Range("K27").Value = 1000
$ws.Range("L27").Value = 900
$ws.Range("M27").Value = -893
$ws.Range("N27").Value = -1114

$ws.Range("H46").Value = 5940.6875
$ws.Range("I46").Value = 1650
$ws.Range("J46").Value = 7891
$ws.Range("K46").Value = 1650
$ws.Range("L46").Value = 7891
$ws.Range("M46").Value = -1462
$ws.Range("N46").Value = -8267

$ws.Range("H132").Value = 2835.6365
$ws.Range("I132").Value = 2761.5
$ws.Range("K132").Value = 8284.5
$ws.Range("M132").Value = -5754.5

$ws.Range("H135").Value = 72333
$ws.Range("J135").Value = 72333
$ws.Range("L135").Value = 72333
$ws.Range("N135").Value = -82473

$ws.Range("H136").Value = 6750
$ws.Range("I136").Value = 6750
$ws.Range("K136").Value = 20250
$ws.Range("M136").Value = -17700

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 188.88889
$ws.Range("I2").Value = 162.5
$ws.Range("K2").Value = 162.5
$ws.Range("M2").Value = -50.5

$ws.Range("H24").Value = 5000
$ws.Range("J24").Value = 5000
$ws.Range("L24").Value = 5000
$ws.Range("N24").Value = -5460

$ws.Range("H107").Value = 27778282
$ws.Range("I107").Value = 27778282
$ws.Range("K107").Value = 83334846
$ws.Range("M107").Value = -83332926

$ws.Range("H113").Value = 731.1177
$ws.Range("I113").Value = 448.27274
$ws.Range("J113").Value = 1249.6666
$ws.Range("K113").Value = 1344.81822
$ws.Range("L113").Value = 3748.9998
$ws.Range("M113").Value = 825.1817799999999
$ws.Range("N113").Value = -8088.9998

$ws.Range("H136").Value = 3953.25
$ws.Range("I136").Value = 2766.8
$ws.Range("J136").Value = 5930.6665
$ws.Range("K136").Value = 8300.400000000001
$ws.Range("L136").Value = 17791.9995
$ws.Range("M136").Value = -5750.400000000001
$ws.Range("N136").Value = -22891.9995

Write-Output "Updated market price data across 49 rows in 8 sheets"
